$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -6
    3  = -3
    4  = -3
    5  = -2
    6  = -2
    7  = -3
    8  = 1
    9  = -6
    10 = 3
    11 = -5
    12 = 1
    14 = -5
    15 = 2
    16 = 1
    17 = -4
    18 = 3
    19 = -1
    20 = 7
    21 = -5
    22 = 1
    23 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
